$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.599.85'
$ws.Range('D3').Value = '1.660.96'
$ws.Range('E3').Value = '  -4.12%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.17'
$ws.Range('E5').Value = '  -2.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.508'
$ws.Range('E6').Value = '  -2.92%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '24.29'
$ws.Range('E8').Value = '  +0.81%  '
$ws.Range('E10').Value = '  -2.61%  '
$ws.Range('E11').Value = '  -1.95%  '
$ws.Range('D12').Value = '1.896.97'
$ws.Range('E12').Value = '  -4.08%  '
$ws.Range('D13').Value = '1.669.50'
$ws.Range('E13').Value = '  -3.67%  '
$ws.Range('E14').Value = '  -2.99%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.567'
$ws.Range('E15').Value = '  +0.85%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.03'
$ws.Range('E16').Value = '  -2.49%  '
$ws.Range('D17').Value = '27.574.71'
$ws.Range('E17').Value = '  -2.51%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '240.87'
$ws.Range('E18').Value = '  -0.74%  '
$ws.Range('D19').Value = '0.0₃0731'
$ws.Range('E19').Value = '  -3.20%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.65'
$ws.Range('E20').Value = '  -3.95%  '
$ws.Range('E21').Value = '  +0.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.50'
$ws.Range('E22').Value = '  -3.39%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.42'
$ws.Range('E23').Value = '  -3.46%  '
$ws.Range('E24').Value = '  -2.39%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.11'
$ws.Range('E25').Value = '  -2.29%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.23'
$ws.Range('E26').Value = '  -4.19%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.36'
$ws.Range('E27').Value = '  -1.82%  '
$ws.Range('E28').Value = '  -0.06%  '
$ws.Range('E29').Value = '  -2.42%  '
$ws.Range('E30').Value = '  +0.46%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0503'
$ws.Range('E31').Value = '  -2.31%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.34'
$ws.Range('E32').Value = '  -2.69%  '
$ws.Range('D33').Value = '1.464.00'
$ws.Range('E33').Value = '  -1.12%  '
$ws.Range('E34').Value = '  -4.65%  '
$ws.Range('E35').Value = '  -4.43%  '
$ws.Range('E36').Value = '  -0.96%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.928'
$ws.Range('E37').Value = '  -4.88%  '
$ws.Range('E38').Value = '  -2.31%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.574'
$ws.Range('E39').Value = '  -4.99%  '
$ws.Range('E40').Value = '  -0.33%  '
$ws.Range('E41').Value = '  -4.83%  '
$ws.Range('E43').Value = '  -4.07%  '
$ws.Range('E44').Value = '  -3.59%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.792'
$ws.Range('E45').Value = '  -0.87%  '
$ws.Range('D46').Value = '1.804.84'
$ws.Range('E46').Value = '  -4.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.71'
$ws.Range('E47').Value = '  -1.59%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '88.81'
$ws.Range('E48').Value = '  -2.41%  '
$ws.Range('E49').Value = '  -5.73%  '
$ws.Range('E50').Value = '  -1.58%  '
$ws.Range('E51').Value = '  -3.52%  '
